$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing D:K to F:M
$ws.Range("D:E").Insert(-4161, 0)

# Copy number formats from column F (the old column D data) into new D:E columns
# so the new quarter columns inherit the same per-row formatting (date vs number).
$ws.Range("F8:F35").Copy()
$ws.Range("D8:E35").PasteSpecial(-4122)
$ws.Range("F41:F102").Copy()
$ws.Range("D41:E102").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# New quarter-end dates for the three "Period Ending" header rows
$ws.Range("D7").Value = (Get-Date -Year 2018 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E7").Value = (Get-Date -Year 2018 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D38").Value = (Get-Date -Year 2018 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E38").Value = (Get-Date -Year 2018 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D80").Value = (Get-Date -Year 2018 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E80").Value = (Get-Date -Year 2018 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0)

# New financial figures for the two newly-added quarter columns (D = latest quarter, E = prior quarter)
$ws.Range("D8").Value = 9600
$ws.Range("E8").Value = 9000
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 2000
$ws.Range("E17").Value = 1700
$ws.Range("D18").Value = 7600
$ws.Range("E18").Value = 7300
$ws.Range("D20").Value = -4200
$ws.Range("E20").Value = -3900
$ws.Range("D21").Value = 4200
$ws.Range("E21").Value = 4100
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 3400
$ws.Range("E23").Value = 3400
$ws.Range("D24").Value = 500
$ws.Range("E24").Value = 600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 2900
$ws.Range("E26").Value = 2900
$ws.Range("D27").Value = 2900
$ws.Range("E27").Value = 2900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 4200
$ws.Range("E32").Value = 3900
$ws.Range("D33").Value = 2800
$ws.Range("E33").Value = 2900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 2800
$ws.Range("E35").Value = 2900
$ws.Range("D41").Value = 16000
$ws.Range("E41").Value = 14300
$ws.Range("D42").Value = 7800
$ws.Range("E42").Value = 7300
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 18600
$ws.Range("E48").Value = 16200
$ws.Range("D49").Value = 200
$ws.Range("E49").Value = 200
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = "NA"
$ws.Range("E52").Value = "NA"
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 981100
$ws.Range("E54").Value = 949900
$ws.Range("D57").Value = 9000
$ws.Range("E57").Value = 8800
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 300
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 887500
$ws.Range("E66").Value = 859800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 64900
$ws.Range("E72").Value = 63100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 93600
$ws.Range("E76").Value = 90100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D81").Value = 2800
$ws.Range("E81").Value = 2900
$ws.Range("D83").Value = 800
$ws.Range("E83").Value = 700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 3400
$ws.Range("E89").Value = 2700
$ws.Range("D91").Value = -1700
$ws.Range("E91").Value = -1000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -29000
$ws.Range("E94").Value = -27700
$ws.Range("D96").Value = -1000
$ws.Range("E96").Value = -600
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 26200
$ws.Range("E100").Value = 23900
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 500
$ws.Range("E102").Value = -1000
